$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new cell at D1 only, shifting D1:H1 right to E1:I1 (row-1-only shift,
# matches "Insert Cells -> Shift cells right" applied to the single cell D1).
$ws.Range("D1").Insert(-4161) | Out-Null   # xlShiftToRight = -4161

# New header text
$ws.Range("D1").Value = "region"
$ws.Range("J1").Value = "ivrticketcode"

# Clear the leftover formatting/content in D2 (style-only cell, no value) so the
# row becomes completely empty and drops out of the saved sheet.
$ws.Range("D2").Clear() | Out-Null

# Update selection to match the recorded cursor position after the edit.
$ws.Range("B3").Select() | Out-Null
